# Apply scheduled market-data refresh values to each sheet (per-cell updates).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5538.684
$ws.Range("I32").Value = 5504.3335
$ws.Range("J32").Value = 5597.5713
$ws.Range("K32").Value = 5504.3335
$ws.Range("L32").Value = 5597.5713
$ws.Range("M32").Value = -5178.3335
$ws.Range("N32").Value = -6249.5713

$ws.Range("H138").Value = 3295.9246
$ws.Range("I138").Value = 668.3570999999999
$ws.Range("J138").Value = 4239.154
$ws.Range("K138").Value = 2005.0713
$ws.Range("L138").Value = 12717.462
$ws.Range("M138").Value = 3134.9287
$ws.Range("N138").Value = -22997.462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16161.25
$ws.Range("I2").Value = 23445.4
$ws.Range("J2").Value = 4021
$ws.Range("K2").Value = 23445.4
$ws.Range("L2").Value = 4021
$ws.Range("M2").Value = -23332.4
$ws.Range("N2").Value = -4247

$ws.Range("H74").Value = 3061.9211
$ws.Range("I74").Value = 1678.6666
$ws.Range("K74").Value = 1678.6666
$ws.Range("M74").Value = -804.6666

$ws.Range("H77").Value = 3061.9211
$ws.Range("I77").Value = 1678.6666
$ws.Range("K77").Value = 8393.333000000001
$ws.Range("M77").Value = -4025.333000000001

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H116").Value = 16161.25
$ws.Range("I116").Value = 23445.4
$ws.Range("J116").Value = 4021
$ws.Range("K116").Value = 23445.4
$ws.Range("L116").Value = 4021
$ws.Range("M116").Value = -21151.4
$ws.Range("N116").Value = -8609

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16161.25
$ws.Range("I3").Value = 23445.4
$ws.Range("J3").Value = 4021
$ws.Range("K3").Value = 23445.4
$ws.Range("L3").Value = 4021
$ws.Range("M3").Value = -23331.4
$ws.Range("N3").Value = -4249

$ws.Range("H20").Value = 2829.2354
$ws.Range("I20").Value = 1239.125
$ws.Range("J20").Value = 4242.6665
$ws.Range("K20").Value = 1239.125
$ws.Range("L20").Value = 4242.6665
$ws.Range("M20").Value = -992.125
$ws.Range("N20").Value = -4736.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 41608.75
$ws.Range("I132").Value = 5184.3335
$ws.Range("J132").Value = 78033.164
$ws.Range("K132").Value = 15553.0005
$ws.Range("L132").Value = 234099.492
$ws.Range("M132").Value = -13023.0005
$ws.Range("N132").Value = -239159.492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 48443040
$ws.Range("I4").Value = 37276804
$ws.Range("K4").Value = 111830412
$ws.Range("M4").Value = -111830300

$ws.Range("H5").Value = 346437.1
$ws.Range("I5").Value = 1594.875
$ws.Range("J5").Value = 770858.3
$ws.Range("K5").Value = 4784.625
$ws.Range("L5").Value = 2312574.9
$ws.Range("M5").Value = -4672.625
$ws.Range("N5").Value = -2312798.9

$ws.Range("H12").Value = 188.5
$ws.Range("J12").Value = 135.53847
$ws.Range("L12").Value = 406.61541
$ws.Range("N12").Value = -752.61541

$ws.Range("H99").Value = 6426.636
$ws.Range("J99").Value = 9281.143
$ws.Range("L99").Value = 27843.429
$ws.Range("N99").Value = -32335.429

$ws.Range("H130").Value = 10821.429
$ws.Range("I130").Value = 4460
$ws.Range("J130").Value = 14355.556
$ws.Range("K130").Value = 13380
$ws.Range("L130").Value = 43066.66800000001
$ws.Range("M130").Value = -8360
$ws.Range("N130").Value = -53106.66800000001

$ws.Range("H131").Value = 6606
$ws.Range("J131").Value = 2105.1667
$ws.Range("L131").Value = 6315.500100000001
$ws.Range("N131").Value = -16395.5001

$ws.Range("H135").Value = 346437.1
$ws.Range("I135").Value = 1594.875
$ws.Range("J135").Value = 770858.3
$ws.Range("K135").Value = 14353.875
$ws.Range("L135").Value = 6937724.7
$ws.Range("M135").Value = -11818.875
$ws.Range("N135").Value = -6942794.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 8904.25
$ws.Range("I36").Value = 7289
$ws.Range("J36").Value = 13750
$ws.Range("K36").Value = 7289
$ws.Range("L36").Value = 13750
$ws.Range("M36").Value = -6804
$ws.Range("N36").Value = -14720

$ws.Range("H92").Value = 13656.125
$ws.Range("I92").Value = 14999
$ws.Range("J92").Value = 13464.286
$ws.Range("K92").Value = 14999
$ws.Range("L92").Value = 13464.286
$ws.Range("M92").Value = -13127
$ws.Range("N92").Value = -17208.286

$ws.Range("H94").Value = 224500000
$ws.Range("J94").Value = 448000000
$ws.Range("L94").Value = 448000000
$ws.Range("N94").Value = -448001352

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H97").Value = 11353.546
$ws.Range("I97").Value = 11353.546
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 11353.546
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -10857.546
$ws.Range("N97").ClearContents()

$ws.Range("H98").Value = 643000000
$ws.Range("J98").Value = 643000000
$ws.Range("L98").Value = 643000000
$ws.Range("N98").Value = -643005990

$ws.Range("H99").Value = 20348.25
$ws.Range("I99").Value = 17131.334
$ws.Range("J99").Value = 29999
$ws.Range("K99").Value = 17131.334
$ws.Range("L99").Value = 29999
$ws.Range("M99").Value = -14885.334
$ws.Range("N99").Value = -34491

$ws.Range("H100").Value = 82245
$ws.Range("J100").Value = 82245
$ws.Range("L100").Value = 82245
$ws.Range("N100").Value = -84409

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H123").Value = 21899.8
$ws.Range("J123").Value = 21899.8
$ws.Range("L123").Value = 21899.8
$ws.Range("N123").Value = -26799.8

$ws.Range("H132").Value = 2688.647
$ws.Range("I132").Value = 2513.8667
$ws.Range("K132").Value = 7541.6001
$ws.Range("M132").Value = -5011.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 38737.383
$ws.Range("I40").Value = 53324.625
$ws.Range("J40").Value = 15397.8
$ws.Range("K40").Value = 53324.625
$ws.Range("L40").Value = 15397.8
$ws.Range("M40").Value = -53188.625
$ws.Range("N40").Value = -15669.8

$ws.Range("H136").Value = 6794.1562
$ws.Range("I136").Value = 6053.5454
$ws.Range("J136").Value = 7182.095
$ws.Range("K136").Value = 18160.6362
$ws.Range("L136").Value = 21546.285
$ws.Range("M136").Value = -15610.6362
$ws.Range("N136").Value = -26646.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H126").Value = 65212.855
$ws.Range("I126").Value = 138833.33
$ws.Range("K126").Value = 416499.99
$ws.Range("M126").Value = -414029.99

$ws.Range("H132").Value = 8514.029
$ws.Range("I132").Value = 9557.654
$ws.Range("K132").Value = 28672.962
$ws.Range("M132").Value = -26142.962

$ws.Range("H133").Value = 42000
$ws.Range("J133").Value = 42000
$ws.Range("L133").Value = 42000
$ws.Range("N133").Value = -52120
